# fix: kolom header "Grade" diganti menjadi "Priority" supaya sesuai dengan
# kolom row data yang sebelumnya meleset; highlight baris yang jadi item
# feedback terkait (#5), geser selection, dan rapikan lebar kolom D baru.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header kolom D: "Grade" -> "Priority"
$ws.Range("D3").Value2 = "Priority"

# 2) Label legend di samping (F8): "Grade:" -> "Priority:"
$ws.Range("F8").Value2 = "Priority:"

# 3) Highlight row 11 (item feedback #5 - "perbaikan row data tergeser dari
#    header column") dengan warna kuning yang sama seperti row 9, menandakan
#    item ini yang sedang dikerjakan/diperbaiki.
$ws.Range("B11:D11").Interior.Color = 65535

# 4) Lebar kolom D baru supaya header "Priority" pas (best-fit).
$ws.Columns.Item(4).ColumnWidth = 9

# 5) Pindahkan selection cursor ke B20.
$ws.Range("B20").Select()
